$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Restore populationsFolder -> populationsCSV
$ws.Range("A7").Value = "populationsCSV"

# Update the current selection to A8
$ws.Range("A8").Select()
